$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Warheads")

# Update "Anti strike craft" (column N) values for the affected weapon rows.
$ws.Range("N2").Value  = 0.075
$ws.Range("N3").Value  = 0.075
$ws.Range("N4").Value  = 0.1
$ws.Range("N6").Value  = 0.075
$ws.Range("N7").Value  = 0.075
$ws.Range("N8").Value  = 0.125
$ws.Range("N10").Value = 0.075
$ws.Range("N11").Value = 0.075
$ws.Range("N12").Value = 0.15
$ws.Range("N20").Value = 0.075
$ws.Range("N21").Value = 0.075
$ws.Range("N22").Value = 0.25
$ws.Range("N26").Value = 0.325
$ws.Range("N40").Value = 0.5
$ws.Range("N44").Value = 0.625

# Update the active cell selection shown when the sheet is saved.
$ws.Range("C2").Select()
